# Apply odds updates to Sheet1 as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("S2").Value = 1.53
$ws.Range("T2").Value = 2.38
$ws.Range("AH2").Value = 9.5
$ws.Range("AN2").Value = 3.75
$ws.Range("AT2").Value = 2.38
$ws.Range("AX2").Value = 26

# Row 3
$ws.Range("G3").Value = 2.75
$ws.Range("I3").Value = 2.8
$ws.Range("J3").Value = 3.4
$ws.Range("L3").Value = 3.5
$ws.Range("AA3").Value = 23
$ws.Range("AH3").Value = 8
$ws.Range("AI3").Value = 13
$ws.Range("AO3").Value = 15
$ws.Range("AX3").Value = 17

# Row 4
$ws.Range("M4").Value = 1.13
$ws.Range("N4").Value = 6

# Row 5
$ws.Range("Q5").Value = 1.93
$ws.Range("R5").Value = 1.93

$wb.Save()
